$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912" ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = "Última actualización: 11:48:04"
$ws1.Cells.Item(3, 1).Value = "Total filas: 209"

$rows1 = @(
    @(96, "08:21:50", "09:01", "23_HERNANDEZ", 40, "LP1912"),
    @(97, "07:20:40", "09:01", "215A_EL PATO", 101, "LP1912"),
    @(109, "07:59:28", "09:22", "16_SANTA ANA", 83, "LP1912"),
    @(110, "07:47:32", "09:22", "17_ROMERO", 95, "LP1912"),
    @(121, "09:38:09", "09:41", "14_ABASTO", 3, "LP1912"),
    @(123, "09:38:09", "09:41", "23_HERNANDEZ", 3, "LP1912"),
    @(175, "11:20:07", "11:35", "23_HERNANDEZ", 15, "LP1912"),
    @(176, "10:26:41", "11:35", "11_ETCHEVERRY", 69, "LP1912"),
    @(180, "11:48:04", "11:49", "16_SANTA ANA", 1, "LP1912"),
    @(181, "10:26:41", "11:51", "215B_EL PATO", 85, "LP1912"),
    @(182, "10:56:30", "11:52", "15_ABASTO", 56, "LP1912"),
    @(183, "11:48:04", "11:53", "16_SANTA ANA", 5, "LP1912"),
    @(184, "10:26:41", "11:59", "225_GOMEZ", 93, "LP1912"),
    @(185, "10:26:41", "12:02", "84_COLONIA URQUIZA-ESC 49", 96, "LP1912"),
    @(186, "11:20:07", "12:05", "23_HERNANDEZ", 45, "LP1912"),
    @(187, "10:26:41", "12:06", "16_P MOR-SANTA ANA", 100, "LP1912"),
    @(188, "10:56:30", "12:06", "14_ABASTO", 70, "LP1912"),
    @(189, "11:20:07", "12:07", "14_ABASTO", 47, "LP1912"),
    @(190, "11:20:07", "12:07", "16_P MOR-SANTA ANA", 47, "LP1912"),
    @(191, "10:56:30", "12:10", "10_OLMOS", 74, "LP1912"),
    @(192, "11:20:07", "12:13", "10_OLMOS", 53, "LP1912"),
    @(193, "10:26:41", "12:14", "17_ROMERO", 108, "LP1912"),
    @(194, "10:26:41", "12:19", "14_ABASTO", 113, "LP1912"),
    @(195, "10:26:41", "12:20", "215A_EL PATO", 114, "LP1912"),
    @(196, "10:56:30", "12:20", "14_ABASTO", 84, "LP1912"),
    @(197, "10:26:41", "12:21", "26_HERNANDEZ", 115, "LP1912"),
    @(198, "11:20:07", "12:21", "14_ABASTO", 61, "LP1912"),
    @(199, "11:20:07", "12:21", "215A_EL PATO", 61, "LP1912"),
    @(200, "11:48:04", "12:35", "23_HERNANDEZ", 47, "LP1912"),
    @(201, "11:48:04", "12:35", "11_ETCHEVERRY", 47, "LP1912"),
    @(202, "10:56:30", "12:36", "27_EL RETIRO", 100, "LP1912"),
    @(203, "11:20:07", "12:37", "27_EL RETIRO", 77, "LP1912"),
    @(204, "10:56:30", "12:38", "17_179 Y 38", 102, "LP1912"),
    @(205, "10:56:30", "12:41", "10_OLMOS", 105, "LP1912"),
    @(206, "11:20:07", "12:49", "11_ETCHEVERRY", 89, "LP1912"),
    @(207, "11:20:07", "13:02", "15_ABASTO", 102, "LP1912"),
    @(208, "11:20:07", "13:07", "16_P MOR-SANTA ANA", 107, "LP1912"),
    @(209, "11:20:07", "13:14", "215D_EL PATO", 114, "LP1912"),
    @(210, "11:48:04", "13:20", "10_OLMOS", 92, "LP1912"),
    @(211, "11:48:04", "13:21", "26_HERNANDEZ", 93, "LP1912"),
    @(212, "11:48:04", "13:27", "14_ABASTO", 99, "LP1912"),
    @(213, "11:48:04", "13:36", "15_ABASTO", 108, "LP1912"),
    @(214, "11:48:04", "13:46", "17_ROMERO", 118, "LP1912"),
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet "LP1912-215" ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = "Última actualización: 11:48:04"

# ---- Sheet "6203-6173" ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 11:48:04"
$ws3.Cells.Item(3, 1).Value = "Total filas: 34"

$rows3 = @(
    @(37, "11:48:04", "12:06", "215A_LA PLATA", 18, "L6173"),
    @(38, "10:56:30", "12:54", "215C_LA PLATA", 118, "L6203"),
    @(39, "11:48:04", "13:31", "215B_LP-P MOR-1 Y 57", 103, "L6173"),
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}
